$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Actualizar" refresh: the timestamp column (D) shifts down one slot as a
# newer snapshot is recorded, with a brand new timestamp written for the
# most recent rows.
$ws.Range("D2:D15").Value = 44265.6424784351
$ws.Range("D16:D29").Value = 44265.61992099537
$ws.Range("D30:D43").Value = 44264.75935748842
